$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes existing rows 6-31 down to 7-32),
# inheriting formatting from the row that was at 6 (now row 7).
$ws.Rows(6).Insert()

# Row height for the newly inserted row (matches target ht="45")
$ws.Rows(6).RowHeight = 45

# Fill in the new entry: Castiblanco & Wilches research project
# Columns are what/when/with/where/why
# (order of assignment matches the order new shared strings were appended
# upstream: "2022 - 2023", then the project blurb, then the student names)
$ws.Range("B6").Value2 = "2022 - 2023"
$ws.Range("E6").Value2 = "Research project: \textit{\href{https://youtu.be/FlZvukFqTcc}{El rol del género en la identificación de la sociosexualidad a partir de las voces} [The role of gender in the identification of sociosexuality from voices]}"
$ws.Range("C6").Value2 = "Maria Camila Wilches \& Johan Sebatián Castiblanco"
$ws.Range("A6").Value2 = "BSc in Psychology"
$ws.Range("D6").Value2 = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"

# Update selection/view to match the author's final state
$ws.Range("B6:E6").Select()
